# "Generate Report for Handoff"
#
# The localization-status report is regenerated: for the two files that were
# still "Ready for handoff" (10ced075-... and its siblings in rows 4-7 of the
# per-locale sheets), the handoff pass re-ran and:
#   - their Priority flips from "low" to "ht"
#   - their "Latest Handoff Datetime" is refreshed to the new handoff run's
#     timestamp (this also updates the Overview sheet's "Latest HO Xliff
#     Generate Date" column for those rows, since it mirrors the same value).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn sheet (rows 4-7: 10ced075, 9b76a4c5, cd8a6af7, f474a370)
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-04 10:35:16"

# de-de sheet (rows 4-7: 10ced075, 9b76a4c5, cd8a6af7, f474a370)
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-04 10:35:20"

# Overview sheet "Latest HO Xliff Generate Date" for the same rows
$wsOverview.Range("G4:G7").Value = "2016-09-04 10:35:20"
